$d = $word.ActiveDocument

# Merge the three "profile" related Class Diagram bullets into one:
#   Visualizza profilo        -> renamed to "Gestisci profilo"
#   Modifica dati personali   -> removed (merged)
#   Cambia foto profilo       -> removed (merged)

$d.Content.Find.Execute("Visualizza profilo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Gestisci profilo", 2) | Out-Null

$toDelete = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Modifica dati personali" -or $t -eq "Cambia foto profilo") {
        $toDelete += $p
    }
}

for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Range.Delete()
}
